$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")  # the "质量监督表" (quality-supervision form) sheet

# Update the text values (stored as plain text / shared strings).
$ws.Range("D2").Value = "芜湖众宇环保建材有限公司"
$ws.Range("B3").Value = "安徽金鹏建设集团股份有限公司"
$ws.Range("D3").Value = "买卖合同纠纷"
$ws.Range("B4").Value = "（2022）皖1103民初26号"

# D4/F4 hold date-looking text ("2021.11.09" / "2022.1.14"). A plain
# Value/Formula assignment gets auto-parsed into a real date serial and
# also stamps a date NumberFormat onto the cell's style, which changes the
# cell's style id. The source file stores these as literal text (General
# style, untouched), so: enter the text with a leading apostrophe to keep
# it literal, then restore the original (General) cell formatting by
# copying it back from a same-style neighbour cell (C4) via Paste Special
# (formats only) - this does not disturb C4 itself.
$ws.Range("D4").Value = "'2021.11.09"
$ws.Range("C4").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F4").Value = "'2022.1.14"
$ws.Range("C4").Copy()
$ws.Range("F4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B5").Value = "调解结案"

$excel.CutCopyMode = $false

# Update the sheet view: drop the scrolled topLeftCell (A12) and change the
# active selection to B5:F5 with B5 as the active cell.
$ws.Activate()
$ws.Range("B5:F5").Select()
